$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell B7 from "DAY CLINIC" to "INTERNAÇÃO" to correct the sector
# mapping for "DIARIA DAY CLINIC" (row 7).
$ws.Range("B7").Value = "INTERNAÇÃO"
